$d = $word.ActiveDocument

# --- Whole-paragraph content replacements ---
# (Paragraph numbering / styles stay put; only each paragraph's inner text is rewritten.
#  NOTE: assigning directly to Paragraphs(n).Range.Text only rewrites the first run and
#  leaves any further runs in that paragraph behind, so we scope a Range that stops one
#  character short of the paragraph mark -- i.e. excludes the pilcrow -- before assigning.)

# Paragraph 6: Objetivos body -> becomes the short "Programa resumido" bullet list
$pr6 = $d.Paragraphs(6).Range
$scoped6 = $d.Range($pr6.Start, $pr6.End - 1)
$scoped6.Text = 'Análise tridimensional de tensões' + [char]11 + 'Critérios de Falha' + [char]11 + 'Fundamentos da Teoria da Elasticidade' + [char]11 + 'Análise Numérica de Tensões e Deformações'

# Paragraph 8: Docentes list -> becomes the old Objetivos paragraph + old Programa body + "A avaliacao..." sentence
$pr8 = $d.Paragraphs(8).Range
$scoped8 = $d.Range($pr8.Start, $pr8.End - 1)
$scoped8.Text = 'Aprofundar conceitos de tensões e deformações em sólidos submetidos aos carregamentos multiaxiais. Apresentar principais critérios de falha local, formulados em termos de diversos parâmetros de tensão ou de deformação. Desenvolver habilidade de formular problemas de contorno da teoria linear de elasticidade para fins de análise de tensões e deformações. Apresentar fundamentos do método dos elementos finitos aplicado aos problemas de contorno da teoria de elasticidade e ensinar principais passos de análise numérica de tensões e deformações.' + [char]11 + 'Análise tridimensional de tensões:' + [char]11 + 'Equações de Transformação no caso Triaxial; Tensões Principais: autovalores e autovetores do tensor de tensão; Invariantes do Estado de Tensão; Tensões Octaédricas.' + [char]11 + 'Critérios de Falha:' + [char]11 + 'Critérios de Fratura para Materiais Frágeis; Critério de Escoamento da Máxima Tensão Cisalhante (Tresca); Critério da Energia de Distorção (Von Mises); Componentes Hidrostático e Desviador do Estado de Tensão. ' + [char]11 + 'Fundamentos da Teoria da Elasticidade: ' + [char]11 + 'Estado de Tensão em um Sólido Contínuo; Estado de deformação: Relações Deformação-Deslocamento (equações cinemáticas), deformação em 3 dimensões e os Invariantes da deformação; Equações Diferenciais de Equilíbrio; Equações de Compatibilidade: Interpretações matemática e física; Princípio de Saint-Venant; Problemas Bidimensionais; Equação de Compatibilidade para o caso bidimensional; Relações Básicas em Coordenadas Polares; Aplicação em Problemas Axissimétricos (tubos de paredes grossas); ' + [char]11 + 'Análise Numérica de Tensões e Deformações:' + [char]11 + 'Diferenças Finitas; Introdução ao Método dos Elementos Finitos; Princípio dos Trabalhos Virtuais, o Problema Unidimensional; Problema Bidimensional; Discretização: Funções de aproximação para elementos triangulares; Emprego de programas computacionais na análise de tensões e deformações pelo Método dos Elementos Finitos' + [char]11 + 'A avaliação será composta por duas provas (P1 e P2).'

# Paragraph 10: Programa resumido body -> becomes the NS/NP1/NP2 criteria text
$pr10 = $d.Paragraphs(10).Range
$scoped10 = $d.Range($pr10.Start, $pr10.End - 1)
$scoped10.Text = 'NS = NP1+NP2; ' + [char]11 + 'NP1: questões da P1 valendo até 4p. no total; ' + [char]11 + 'NP2: questões da P2 valendo até 6 p. no total.'

# Paragraph 12: Programa body -> becomes the "A recuperacao..." sentence
$pr12 = $d.Paragraphs(12).Range
$scoped12 = $d.Range($pr12.Start, $pr12.End - 1)
$scoped12.Text = 'A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2.'

# Paragraph 16: Bibliografia body -> becomes the Viktor Pastoukhov line
$pr16 = $d.Paragraphs(16).Range
$scoped16 = $d.Range($pr16.Start, $pr16.End - 1)
$scoped16.Text = '7797767 - Viktor Pastoukhov'

# --- Paragraph 14 (Avaliacao bullet): bold labels (Metodo: / Criterio: / Norma de recuperacao:)
# stay exactly where they are; only the plain-text runs between them are swapped out, via
# Find/Replace scoped to the paragraph's own Range so the bold runs are never touched. ---
$p14 = $d.Paragraphs(14).Range
$p14.Find.Execute('A avaliação será composta por duas provas (P1 e P2).' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.' + [char]11 + [char]11 + '2. M.H. SADD. Elasticity: Theory, Applications and Numerics. Amsterdam: Elsevier, 2005, 461p.' + [char]11 + [char]11 + '3. R.R. CRAIG,Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p. ' + [char]11 + [char]11 + '4. A.C. UGURAL, S.K. FENSTER. Advanced Strength and Applied Elasticity. New Jersey: Prentice Hall. 4th Ed., 2003, 544p.' + [char]11 + [char]11 + '5. S.P. TIMOSHENKO, J.N. GOODIER. Teoria da Elasticidade. Rio de janeiro: Guanabara Dois. 3a Ed., 1980, 545p.' + [char]11 + [char]11 + '6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p. ' + [char]11 + [char]11 + '7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.' + [char]11 + [char]11 + '8. T.M. ATANACKOVIC, A. GURAN. Theory of Elasticity for Scientists and Engineers. New York: Springer Science+Business, 2000, 374p.' + [char]11, 2) | Out-Null
$p14.Find.Execute('NS = NP1+NP2; ' + [char]11 + 'NP1: questões da P1 valendo até 4p. no total; ' + [char]11 + 'NP2: questões da P2 valendo até 6 p. no total.' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '471420 - Carlos Antonio Reis Pereira Baptista' + [char]11, 2) | Out-Null
$p14.Find.Execute('A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2.', $true, $false, $false, $false, $false, $true, 1, $false, '3480026 - João Paulo Pascon', 2) | Out-Null
